$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 - Mazda2 Sedan: mark VERIFICAR = TRUE, bump all 4 price tiers
# ---------------------------------------------------------------------
$ws.Range("A2").Value = $true
$ws.Range("H2").Value = "300,900"
$ws.Range("M2").Value = "330,900"
$ws.Range("R2").Value = "370,900"
$ws.Range("W2").Value = "398,900"

# ---------------------------------------------------------------------
# Row 3 - Mazda2 Hatchback: mark VERIFICAR = TRUE, bump all 3 price tiers
# ---------------------------------------------------------------------
$ws.Range("A3").Value = $true
$ws.Range("H3").Value = "330,900"
$ws.Range("M3").Value = "370,900"
$ws.Range("R3").Value = "398,900"

# ---------------------------------------------------------------------
# Row 11 - Mazda CX-90: no longer flagged for VERIFICAR
# ---------------------------------------------------------------------
$ws.Range("A11").Value = $false

# ---------------------------------------------------------------------
# Insert a new row at 12 for the Mazda BT-50 pickup (new PICKUPS
# category). Inserting above row 12 shifts the two MX-5 rows down to
# 13/14 and copies formatting down from row 11.
# ---------------------------------------------------------------------
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = $false
$ws.Range("B12").Value = "mazda-bt-50-roja-flydown-inclinado-v2"
$ws.Range("C12").Value = "PICKUPS"
$ws.Range("D12").Value = "MAZDA BT-50"
$ws.Range("E12").Value = "N/A"
$ws.Range("F12").Value = 2025
$ws.Range("G12").Value = "Signature"
$ws.Range("H12").Value = "829,900"
$ws.Range("I12").Value = "188"
$ws.Range("J12").Value = "332"
$ws.Range("K12").Value = "3.0L"

# ---------------------------------------------------------------------
# Row 13 (previously row 12) - Mazda MX-5: refreshed render, 2025 model
# year, new price
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "mx-5-carro-rojo-deportivo-inclinado-v1-1"
$ws.Range("F13").Value = 2025
$ws.Range("H13").Value = "529,900"

# ---------------------------------------------------------------------
# Row 14 (previously row 13) - Mazda MX-5 RF: 2025 model year, new price
# ---------------------------------------------------------------------
$ws.Range("F14").Value = 2025
$ws.Range("H14").Value = "629,900"

# ---------------------------------------------------------------------
# Insert a new row at 15 for the Mazda MX-5 35th Anniversary edition
# ---------------------------------------------------------------------
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = $false
$ws.Range("B15").Value = "mazda-mx-5-35-aniversario-flydown-perspectiva-v1"
$ws.Range("C15").Value = "ROADSTERS"
$ws.Range("D15").Value = "MAZDA MX-5 35° ANIVERSARIO"
$ws.Range("E15").Value = "N/A"
$ws.Range("F15").Value = 2025
$ws.Range("G15").Value = "35° Aniversario"
$ws.Range("H15").Value = "639,900"
$ws.Range("I15").Value = "181"
$ws.Range("J15").Value = "151"
$ws.Range("K15").Value = "2.0L"

# ---------------------------------------------------------------------
# Match the selection left behind by the author at the end of their
# editing session
# ---------------------------------------------------------------------
$ws.Range("A13").Select()
